# Fix read database of TeacherView
# - "- master.xml ..." -> "- teacher.xml ..." (renaming the data file reference)
# - add a new line right after it describing lanhdao.xml (Phong ban / leadership data)

$d = $word.ActiveDocument

# --- Step 1: insert the new "lanhdao.xml" paragraph right after the
# paragraph that talks about master.xml / the teacher-data file, BEFORE we
# touch that paragraph's runs below (keeps the new paragraph's own runs free
# of any direct character formatting picked up from later edits).
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*master.xml*") {
        $p.Range.InsertParagraphAfter()
        $newPara = $p.Next()
        $newPara.Range.Text = "- lanhdao.xml – Dữ liệu lãnh đạo trung tâm, quản lý theo mã Lãnh đạo"
        break
    }
}

# --- Step 2: rename master -> teacher inside the original line, splitting
# that run the same way Word would when you select just the word "master"
# and retype it (the surrounding "- " and ".xml – ..." stay their own runs).
$rng = $d.Content
$found = $rng.Find.Execute("master", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)
if ($found) {
    $rng.Text = "teacher"
    # Toggling a character attribute on exactly this sub-range forces Word
    # to materialize it as its own run instead of re-merging with neighbors.
    $rng.Font.Bold = 1
    $rng.Font.Bold = 0
}
